# Auto-generated Excel COM-interop edit script
# Applies the crypto price/volume update described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.283.39"
$ws.Range("E2").Value = "  +2.11%  "
$ws.Range("D3").Value = "3.630.23"
$ws.Range("E3").Value = "  +1.14%  "
$ws.Range("E4").Value = "  -0.35%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "196.23"
$ws.Range("E5").Value = "  +7.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "578.88"
$ws.Range("E6").Value = "  -0.76%  "
$ws.Range("D7").Value = "3.625.58"
$ws.Range("E7").Value = "  +1.06%  "
$ws.Range("E8").Value = "  +3.07%  "
$ws.Range("E9").Value = "  +0.32%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.679"
$ws.Range("E10").Value = "  +1.62%  "
$ws.Range("E11").Value = "  +8.15%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "55.82"
$ws.Range("E12").Value = "  +5.32%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000296"
$ws.Range("E13").Value = "  +20.51%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.14"
$ws.Range("E14").Value = "  +3.21%  "
$ws.Range("D15").Value = "4.207.16"
$ws.Range("E15").Value = "  +0.91%  "
$ws.Range("D16").Value = "3.622.89"
$ws.Range("E16").Value = "  +0.79%  "
$ws.Range("E17").Value = "  +0.49%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.56"
$ws.Range("E18").Value = "  +3.86%  "
$ws.Range("D19").Value = "68.205.24"
$ws.Range("E19").Value = "  +2.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.54"
$ws.Range("E20").Value = "  +1.76%  "
$ws.Range("E21").Value = "  +2.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "404.23"
$ws.Range("E22").Value = "  +3.50%  "
$ws.Range("B23").Value = "RenderToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.71"
$ws.Range("E23").Value = "  +24.78%  "
$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.23"
$ws.Range("E24").Value = "  -1.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.16"
$ws.Range("E25").Value = "  +1.90%  "
$ws.Range("E26").Value = "  +5.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.63"
$ws.Range("E27").Value = "  +4.45%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.87"
$ws.Range("E28").Value = "  +8.75%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.14"
$ws.Range("E30").Value = "  +21.74%  "
$ws.Range("E31").Value = "  +3.89%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "31.75"
$ws.Range("E32").Value = "  +3.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "691.76"
$ws.Range("E33").Value = "  +19.85%  "
$ws.Range("E34").Value = "  +4.10%  "
$ws.Range("E35").Value = "  +6.99%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "64.89"
$ws.Range("E36").Value = "  -1.15%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "42.74"
$ws.Range("E37").Value = "  +3.93%  "
$ws.Range("E38").Value = "  +13.29%  "
$ws.Range("E39").Value = "  +0.25%  "
$ws.Range("E40").Value = "  +12.12%  "
$ws.Range("E41").Value = "  +23.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.15"
$ws.Range("E42").Value = "  +15.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.136"
$ws.Range("E43").Value = "  +3.34%  "
$ws.Range("D44").Value = "3.146.68"
$ws.Range("E44").Value = "  +18.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.999"
$ws.Range("E45").Value = "  -0.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.94"
$ws.Range("E46").Value = "  +30.59%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0425"
$ws.Range("E47").Value = "  +4.80%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.133"
$ws.Range("E48").Value = "  +2.09%  "
$ws.Range("E49").Value = "  +6.49%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.11"
$ws.Range("E50").Value = "  +2.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "142.52"
$ws.Range("E51").Value = "  +2.21%  "
